$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @(
  "https://unitedwayshr.org/staging/",
  "https://unitedwayshr.org/staging/what-we-do/",
  "https://unitedwayshr.org/staging/what-we-do/united-for-children/",
  "https://unitedwayshr.org/staging/what-we-do/mission-united/",
  "https://unitedwayshr.org/staging/what-we-do/developmental-screening/",
  "https://unitedwayshr.org/staging/what-we-do/project-inclusion/",
  "https://unitedwayshr.org/staging/get-involved/give/",
  "https://unitedwayshr.org/staging/get-involved/",
  "https://unitedwayshr.org/staging/about/",
  "https://unitedwayshr.org/staging/about/events/"
)
for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = 3 + $i
    $ws.Range("A$row").Value = $urls[$i]
}

$ws.Range("A12:I12").Copy()
$ws.Range("A13:I15").PasteSpecial(-4122)  # xlPasteFormats

$newUrls = @(
  "https://unitedwayshr.org/staging/about/our-team/",
  "https://unitedwayshr.org/staging/give/",
  "https://unitedwayshr.org/staging/toolkit/"
)
for ($i = 0; $i -lt $newUrls.Length; $i++) {
    $row = 13 + $i
    $ws.Range("A$row").Value = $newUrls[$i]
}

$chk9 = [string]$ws.Range("A9").Value2
$chk10 = [string]$ws.Range("A10").Value2
Write-Output "after newUrls: A9=$chk9 A10=$chk10"

$ws.Range("A1:A20").Hyperlinks.Delete()

$chk9b = [string]$ws.Range("A9").Value2
$chk10b = [string]$ws.Range("A10").Value2
Write-Output "after hyperlink delete: A9=$chk9b A10=$chk10b"
